# Update the "想去人数" (wanted-to-go count) column F values across the
# 展览 (sheet1), 演出 (sheet2) and 全部类型 (sheet4) worksheets to reflect the
# newly generated gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

# --- 展览 sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6926
$ws1.Range("F3").Value = 96
$ws1.Range("F4").Value = 0
$ws1.Range("F6").Value = 157
$ws1.Range("F7").Value = 6739
$ws1.Range("F10").Value = 1290
$ws1.Range("F11").Value = 19
$ws1.Range("F13").Value = 405
$ws1.Range("F14").Value = 144
$ws1.Range("F15").Value = 17
$ws1.Range("F16").Value = 405
$ws1.Range("F17").Value = 0
$ws1.Range("F19").Value = 11
$ws1.Range("F20").Value = 5113
$ws1.Range("F21").Value = 112
$ws1.Range("F23").Value = 0
$ws1.Range("F25").Value = 215

# --- 演出 sheet ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 46

# --- 全部类型 sheet ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 0
$ws4.Range("F3").Value = 96
$ws4.Range("F4").Value = 0
$ws4.Range("F5").Value = 452
$ws4.Range("F6").Value = 157
$ws4.Range("F7").Value = 0
$ws4.Range("F8").Value = 67
$ws4.Range("F9").Value = 0
$ws4.Range("F11").Value = 0
$ws4.Range("F14").Value = 144
$ws4.Range("F15").Value = 17
$ws4.Range("F16").Value = 405
$ws4.Range("F17").Value = 0
$ws4.Range("F18").Value = 0
$ws4.Range("F23").Value = 0
$ws4.Range("F24").Value = 147
$ws4.Range("F25").Value = 0
$ws4.Range("F27").Value = 215
